# Update (Removed Auto Arima)
# Applies updated forecast figures to the "Forecast Comparison" sheet and
# the recomputed summary metrics on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Forecast Comparison
# Columns: C = Prophet Forecast, D = Amazon Mean Forecast,
#          E = Amazon P70 Forecast, F = Amazon P80 Forecast,
#          G = Amazon P90 Forecast
# ---------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$forecastRows = @{
    2  = @(28, 50, 60, 68, 81)
    3  = @(30, 44, 53, 62, 76)
    4  = @(29, 38, 45, 52, 64)
    5  = @(28, 37, 44, 52, 63)
    6  = @(29, 37, 44, 52, 64)
    7  = @(28, 37, 44, 52, 64)
    8  = @(25, 37, 45, 54, 67)
    9  = @(21, 37, 46, 56, 73)
    10 = @(21, 36, 44, 52, 65)
    11 = @(23, 37, 45, 55, 72)
    12 = @(26, 37, 45, 56, 73)
    13 = @(28, 40, 49, 60, 79)
    14 = @(28, 38, 47, 59, 77)
    15 = @(30, 37, 45, 58, 78)
    16 = @(33, 37, 46, 58, 79)
    17 = @(36, 36, 44, 56, 75)
}

foreach ($rowNum in $forecastRows.Keys) {
    $vals = $forecastRows[$rowNum]
    $wsForecast.Range("C$rowNum").Value = $vals[0]
    $wsForecast.Range("D$rowNum").Value = $vals[1]
    $wsForecast.Range("E$rowNum").Value = $vals[2]
    $wsForecast.Range("F$rowNum").Value = $vals[3]
    $wsForecast.Range("G$rowNum").Value = $vals[4]
}

# ---------------------------------------------------------------------
# Sheet: Summary
# Recomputed forecast totals and min-forecast week after removing the
# Auto ARIMA model from the ensemble.
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B9").Value = "'443"
$wsSummary.Range("B10").Value = "'218"
$wsSummary.Range("B11").Value = "'115"
$wsSummary.Range("B15").Value = "'2025-01-26"
